# Update for May metrics
# Fills in March and April actuals for both the Engines (Sheet1) and
# Drivers (Sheet2) YTD reports, refreshes the running totals, and leaves
# May (not yet closed) blank.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Sheet1 ("Engines (Outliers Removed)")
# ---------------------------------------------------------------------

# February (row 4) keeps its values; only re-stamp so the Qty/Qty Pad/Total
# columns share one plain-number style and Cost switches from the 2-decimal
# dollar format to a whole-number format.
$ws1.Range("B4:D4").ClearFormats()
$ws1.Range("E4").Value = 3881.8833
$ws1.Range("E4").NumberFormat = "#,##0"

# March (row 5) - newly closed-out actuals
$ws1.Range("B5").Value = 128
$ws1.Range("B5").ClearFormats()
$ws1.Range("E5").Value = 4402.3
$ws1.Range("E5").NumberFormat = "#,##0"

# April (row 6) - newly closed-out actuals
$ws1.Range("B6").Value = 85
$ws1.Range("B6").ClearFormats()
$ws1.Range("E6").Value = 2309.9286000000006
$ws1.Range("E6").NumberFormat = "`"$`"#,##0"

# May (row 7) not closed yet - clears the inherited "$" format on the Cost cell
$ws1.Range("E7").ClearFormats()

# Running totals recalc automatically off the existing SUM formulas
$ws1.Range("B15").Formula = "=SUM(B3:B14)"
$ws1.Range("E15").Formula = "=SUM(E3:E14)"

# ---------------------------------------------------------------------
# Sheet2 ("Drivers (Outliers Removed)")
# ---------------------------------------------------------------------

$ws2.Range("B4:D4").NumberFormat = "General"

# March (row 5)
$ws2.Range("B5").Value = 303
$ws2.Range("E5").Value = 6420
$ws2.Range("E5").NumberFormat = "`"$`"#,##0"

# April (row 6)
$ws2.Range("B6").Value = 102
$ws2.Range("B6").NumberFormat = "General"
$ws2.Range("E6").Value = 1441.0014999999999
$ws2.Range("E6").NumberFormat = "`"$`"#,##0"
$ws2.Range("E6").HorizontalAlignment = -4108

$ws2.Range("B15").Formula = "=SUM(B3:B14)"
$ws2.Range("E15").Formula = "=SUM(E3:E14)"

# ---------------------------------------------------------------------
# Restore the selections Excel left behind on save
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("E4:E5").Select()

$ws2.Activate()
$ws2.Range("E7").Select()
